$wb = $excel.ActiveWorkbook

# --- Rename "Sheet3" to "Stock Screen" and make it the active tab ---
$wsStock = $wb.Worksheets.Item("Sheet3")
$wsStock.Name = "Stock Screen"

# Activate the Stock Screen sheet (moves tabSelected / activeTab) and
# move its selection to A5 (matching the recorded view state).
$wsStock.Activate()
$wsStock.Range("A5").Select()

# --- Row 9: replace the "use for CUSIP" comment with a "Cusip" UI Field ---
$wsStock.Cells.Item(9, 4).Value = "Cusip"
$wsStock.Cells.Item(9, 5).Value = ""

# --- Row 33: highlight the "WMS_PUT_CALL_STATUS" column name in red ---
$wsStock.Cells.Item(33, 3).Font.Color = 255

Write-Host "done"
